$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.9523829999999999
$ws.Range("H2").Value = 2.857149
$ws.Range("I2").Value = 0.04618630532204829
$ws.Range("J2").Value = 0.0649320951835329
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.9703876666666668
$ws.Range("N2").Value = 2.911163
$ws.Range("O2").Value = 0.03945299285965207
$ws.Range("P2").Value = 0.04754668824173519
$ws.Range("Q2").Value = 0.924180717143
$ws.Range("R2").Value = 8.317626454287
$ws.Range("S2").Value = 0.001822187974084482
$ws.Range("T2").Value = 0.003087306086574114

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.9523829999999999
$ws.Range("H3").Value = 2.857149
$ws.Range("I3").Value = 0.04618630532204829
$ws.Range("J3").Value = 0.0649320951835329
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.041192666666667
$ws.Range("N3").Value = 3.123578
$ws.Range("O3").Value = 0.04233170747586662
$ws.Range("P3").Value = 0.05101596487889641
$ws.Range("Q3").Value = 0.9916141954579999
$ws.Range("R3").Value = 8.924527759122
$ws.Range("S3").Value = 0.00195514516628401
$ws.Range("T3").Value = 0.003312573487396274

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.9523829999999999
$ws.Range("H4").Value = 2.857149
$ws.Range("I4").Value = 0.04618630532204829
$ws.Range("J4").Value = 0.0649320951835329
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.170211333333333
$ws.Range("N4").Value = 15.510634
$ws.Range("O4").Value = 0.2102049704707969
$ws.Range("P4").Value = 0.2533280614069559
$ws.Range("Q4").Value = 4.924021380274
$ws.Range("R4").Value = 44.316192422466
$ws.Range("S4").Value = 0.009708590946376373
$ws.Range("T4").Value = 0.01644912179593633

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.9523829999999999
$ws.Range("H5").Value = 2.857149
$ws.Range("I5").Value = 0.04618630532204829
$ws.Range("J5").Value = 0.0649320951835329
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.853575333333334
$ws.Range("N5").Value = 14.560726
$ws.Range("O5").Value = 0.1973315197085667
$ws.Range("P5").Value = 0.2378136503161547
$ws.Range("Q5").Value = 4.622462636686
$ws.Range("R5").Value = 41.602163730174
$ws.Range("S5").Value = 0.009114013818923654
$ws.Range("T5").Value = 0.01544173857827197

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.9523829999999999
$ws.Range("H6").Value = 2.857149
$ws.Range("I6").Value = 0.04618630532204829
$ws.Range("J6").Value = 0.0649320951835329
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 12.56068
$ws.Range("N6").Value = 25.12136
$ws.Range("O6").Value = 0.5106788094851177
$ws.Range("P6").Value = 0.4102956351562577
$ws.Range("Q6").Value = 11.96257810044
$ws.Range("R6").Value = 71.77546860263999
$ws.Range("S6").Value = 0.02358636741637977
$ws.Range("T6").Value = 0.02664135523535422

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.808798
$ws.Range("H7").Value = 5.426394
$ws.Range("I7").Value = 0.08771859293363103
$ws.Range("J7").Value = 0.1233212309583266
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.9703876666666668
$ws.Range("N7").Value = 2.911163
$ws.Range("O7").Value = 0.03945299285965207
$ws.Range("P7").Value = 0.04754668824173519
$ws.Range("Q7").Value = 1.755235270691334
$ws.Range("R7").Value = 15.797117436222
$ws.Range("S7").Value = 0.003460761020669272
$ws.Range("T7").Value = 0.005863516121962577

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.808798
$ws.Range("H8").Value = 5.426394
$ws.Range("I8").Value = 0.08771859293363103
$ws.Range("J8").Value = 0.1233212309583266
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.041192666666667
$ws.Range("N8").Value = 3.123578
$ws.Range("O8").Value = 0.04233170747586662
$ws.Range("P8").Value = 0.05101596487889641
$ws.Range("Q8").Value = 1.883307213081334
$ws.Range("R8").Value = 16.949764917732
$ws.Range("S8").Value = 0.00371327781626109
$ws.Range("T8").Value = 0.006291351587392263

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.808798
$ws.Range("H9").Value = 5.426394
$ws.Range("I9").Value = 0.08771859293363103
$ws.Range("J9").Value = 0.1233212309583266
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.170211333333333
$ws.Range("N9").Value = 15.510634
$ws.Range("O9").Value = 0.2102049704707969
$ws.Range("P9").Value = 0.2533280614069559
$ws.Range("Q9").Value = 9.351867919310667
$ws.Range("R9").Value = 84.166811273796
$ws.Range("S9").Value = 0.01843888423735377
$ws.Range("T9").Value = 0.03124072836899235

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.808798
$ws.Range("H10").Value = 5.426394
$ws.Range("I10").Value = 0.08771859293363103
$ws.Range("J10").Value = 0.1233212309583266
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.853575333333334
$ws.Range("N10").Value = 14.560726
$ws.Range("O10").Value = 0.1973315197085667
$ws.Range("P10").Value = 0.2378136503161547
$ws.Range("Q10").Value = 8.779137355782668
$ws.Range("R10").Value = 79.01223620204401
$ws.Range("S10").Value = 0.01730964325029056
$ws.Range("T10").Value = 0.02932747209568123

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.808798
$ws.Range("H11").Value = 5.426394
$ws.Range("I11").Value = 0.08771859293363103
$ws.Range("J11").Value = 0.1233212309583266
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 12.56068
$ws.Range("N11").Value = 25.12136
$ws.Range("O11").Value = 0.5106788094851177
$ws.Range("P11").Value = 0.4102956351562577
$ws.Range("Q11").Value = 22.71973286264
$ws.Range("R11").Value = 136.31839717584
$ws.Range("S11").Value = 0.04479602660905635
$ws.Range("T11").Value = 0.05059816278429816

$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 17.859282
$ws.Range("H12").Value = 35.718564
$ws.Range("I12").Value = 0.8660951017443207
$ws.Range("J12").Value = 0.8117466738581405
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.9703876666666668
$ws.Range("N12").Value = 2.911163
$ws.Range("O12").Value = 0.03945299285965207
$ws.Range("P12").Value = 0.04754668824173519
$ws.Range("Q12").Value = 17.330426988322
$ws.Range("R12").Value = 103.982561929932
$ws.Range("S12").Value = 0.03417004386489832
$ws.Range("T12").Value = 0.0385958660331985

$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 17.859282
$ws.Range("H13").Value = 35.718564
$ws.Range("I13").Value = 0.8660951017443207
$ws.Range("J13").Value = 0.8117466738581405
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.041192666666667
$ws.Range("N13").Value = 3.123578
$ws.Range("O13").Value = 0.04233170747586662
$ws.Range("P13").Value = 0.05101596487889641
$ws.Range("Q13").Value = 18.594953450332
$ws.Range("R13").Value = 111.569720701992
$ws.Range("S13").Value = 0.03666328449332152
$ws.Range("T13").Value = 0.04141203980410788

$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 17.859282
$ws.Range("H14").Value = 35.718564
$ws.Range("I14").Value = 0.8660951017443207
$ws.Range("J14").Value = 0.8117466738581405
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 5.170211333333333
$ws.Range("N14").Value = 15.510634
$ws.Range("O14").Value = 0.2102049704707969
$ws.Range("P14").Value = 0.2533280614069559
$ws.Range("Q14").Value = 92.336262201596
$ws.Range("R14").Value = 554.017573209576
$ws.Range("S14").Value = 0.1820574952870668
$ws.Range("T14").Value = 0.2056382112420272

$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 17.859282
$ws.Range("H15").Value = 35.718564
$ws.Range("I15").Value = 0.8660951017443207
$ws.Range("J15").Value = 0.8117466738581405
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.853575333333334
$ws.Range("N15").Value = 14.560726
$ws.Range("O15").Value = 0.1973315197085667
$ws.Range("P15").Value = 0.2378136503161547
$ws.Range("Q15").Value = 86.681370586244
$ws.Range("R15").Value = 520.088223517464
$ws.Range("S15").Value = 0.1709078626393525
$ws.Range("T15").Value = 0.1930444396422015

$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 17.859282
$ws.Range("H16").Value = 35.718564
$ws.Range("I16").Value = 0.8660951017443207
$ws.Range("J16").Value = 0.8117466738581405
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 12.56068
$ws.Range("N16").Value = 25.12136
$ws.Range("O16").Value = 0.5106788094851177
$ws.Range("P16").Value = 0.4102956351562577
$ws.Range("Q16").Value = 224.32472623176
$ws.Range("R16").Value = 897.29890492704
$ws.Range("S16").Value = 0.4422964154596816
$ws.Range("T16").Value = 0.3330561171366053
